# Weekly fruit/vegetable price update: a new price record (week) is
# inserted into the data table at row 378, pushing the existing rows
# 378-471 down to 379-472 (the table otherwise keeps growing by one row
# at the bottom of the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 378, shifting
# everything from 378 downward (including the last row, 471) down by
# one row.
$ws.Rows("378:378").Insert()

# Populate the newly inserted row 378 with the new weekly record.
# Columns that are identical across every row in this table
# (market/category metadata) are copied from the row immediately below
# (old row 378, now row 379). Note: `.Value` reads back a placeholder
# in this host, so use `.Value2` to read existing cell contents.
$ws.Cells.Item(378, 1).Value = $ws.Cells.Item(379, 1).Value2    # A: Mercado ID
$ws.Cells.Item(378, 2).Value = $ws.Cells.Item(379, 2).Value2    # B: Mercado
$ws.Cells.Item(378, 3).Value = $ws.Cells.Item(379, 3).Value2    # C: Region
$ws.Cells.Item(378, 4).Value = 44543                            # D: Fecha
$ws.Cells.Item(378, 5).Value = $ws.Cells.Item(379, 5).Value2    # E: Codreg
$ws.Cells.Item(378, 6).Value = $ws.Cells.Item(379, 6).Value2    # F: Categoria ID
$ws.Cells.Item(378, 7).Value = $ws.Cells.Item(379, 7).Value2    # G: Categoria
$ws.Cells.Item(378, 8).Value = $ws.Cells.Item(379, 8).Value2    # H: Variedad
$ws.Cells.Item(378, 9).Value = "Extra"                          # I: Calidad
$ws.Cells.Item(378, 10).Value = 250                             # J: Volumen
$ws.Cells.Item(378, 11).Value = 20000                           # K: Precio minimo
$ws.Cells.Item(378, 12).Value = 20000                           # L: Precio maximo
$ws.Cells.Item(378, 13).Value = 20000                           # M: Precio promedio ponderado
$ws.Cells.Item(378, 14).Value = $ws.Cells.Item(379, 14).Value2  # N: Unidad de comercializacion
$ws.Cells.Item(378, 15).Value = "Provincia de Quillota"         # O: Origen
$ws.Cells.Item(378, 16).Value = 1111                            # P: Precio $/Kg
$ws.Cells.Item(378, 17).Value = $ws.Cells.Item(379, 17).Value2  # Q: Kg o Unidades
$ws.Cells.Item(378, 18).Value = $ws.Cells.Item(379, 18).Value2  # R: Clasificacion
